$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report header description (A1)
$ws.Range("A1").Value = "Description unknown, completed 05/17/2023 10:26:12 EDT, by WPJTOWN1.The search returned: 10 events."

# Update the data table (rows 3-12) with corrected/re-mapped trace events
$ws.Cells.Item(3,1).Value = "FURX"
$ws.Cells.Item(3,2).Value = 855167
$ws.Cells.Item(3,3).Value = "DENVER"
$ws.Cells.Item(3,4).Value = "CO"
$ws.Cells.Item(3,5).Value = 5
$ws.Cells.Item(3,6).Value = 15
$ws.Cells.Item(3,7).Value = 2233
$ws.Cells.Item(3,8).Value = "Arrive In-Transit"
$ws.Cells.Item(3,9).Value = "HKCKDE"
$ws.Cells.Item(3,10).Value = "LOVELAND"
$ws.Cells.Item(3,11).Value = "CO"
$ws.Cells.Item(3,12).Value = 230708
$ws.Cells.Item(3,13).Value = 60300
$ws.Cells.Item(3,14).Value = 170408
$ws.Cells.Item(3,15).Value = "FURX855167"
$ws.Cells.Item(4,1).Value = "CEFX"
$ws.Cells.Item(4,2).Value = 360837
$ws.Cells.Item(4,3).Value = "JOHNSTOWN"
$ws.Cells.Item(4,4).Value = "CO"
$ws.Cells.Item(4,5).Value = 5
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1334
$ws.Cells.Item(4,8).Value = "Placed Actual"
$ws.Cells.Item(4,10).Value = "LOVELAND"
$ws.Cells.Item(4,11).Value = "CO"
$ws.Cells.Item(4,12).Value = 283534
$ws.Cells.Item(4,13).Value = 68500
$ws.Cells.Item(4,14).Value = 215034
$ws.Cells.Item(4,15).Value = "CEFX360837"
$ws.Cells.Item(5,1).Value = "CRDX"
$ws.Cells.Item(5,2).Value = 15033
$ws.Cells.Item(5,3).Value = "JOHNSTOWN"
$ws.Cells.Item(5,4).Value = "CO"
$ws.Cells.Item(5,5).Value = 5
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1334
$ws.Cells.Item(5,8).Value = "Placed Actual"
$ws.Cells.Item(5,10).Value = "LOVELAND"
$ws.Cells.Item(5,11).Value = "CO"
$ws.Cells.Item(5,12).Value = 278305
$ws.Cells.Item(5,13).Value = 69400
$ws.Cells.Item(5,14).Value = 208905
$ws.Cells.Item(5,15).Value = "CRDX15033"
$ws.Cells.Item(6,1).Value = "CRDX"
$ws.Cells.Item(6,2).Value = 10122
$ws.Cells.Item(6,3).Value = "KANSAS CITY"
$ws.Cells.Item(6,4).Value = "KS"
$ws.Cells.Item(6,5).Value = 5
$ws.Cells.Item(6,6).Value = 16
$ws.Cells.Item(6,7).Value = 1809
$ws.Cells.Item(6,8).Value = "Arrive In-Transit"
$ws.Cells.Item(6,9).Value = "HNTWKC"
$ws.Cells.Item(6,10).Value = "LOVELAND"
$ws.Cells.Item(6,11).Value = "CO"
$ws.Cells.Item(6,12).Value = 253156
$ws.Cells.Item(6,13).Value = 71800
$ws.Cells.Item(6,14).Value = 181356
$ws.Cells.Item(6,15).Value = "CRDX10122"
$ws.Cells.Item(7,1).Value = "AEX"
$ws.Cells.Item(7,2).Value = 9203
$ws.Cells.Item(7,3).Value = "LA CYGNE"
$ws.Cells.Item(7,4).Value = "KS"
$ws.Cells.Item(7,5).Value = 5
$ws.Cells.Item(7,6).Value = 17
$ws.Cells.Item(7,7).Value = 602
$ws.Cells.Item(7,8).Value = "Departure"
$ws.Cells.Item(7,9).Value = "HTULKC"
$ws.Cells.Item(7,10).Value = "LOVELAND"
$ws.Cells.Item(7,11).Value = "CO"
$ws.Cells.Item(7,12).Value = 278800
$ws.Cells.Item(7,13).Value = 67900
$ws.Cells.Item(7,14).Value = 210900
$ws.Cells.Item(7,15).Value = "AEX9203"
$ws.Cells.Item(8,1).Value = "CRDX"
$ws.Cells.Item(8,2).Value = 15088
$ws.Cells.Item(8,3).Value = "LONGMONT"
$ws.Cells.Item(8,4).Value = "CO"
$ws.Cells.Item(8,5).Value = 5
$ws.Cells.Item(8,6).Value = 16
$ws.Cells.Item(8,7).Value = 1453
$ws.Cells.Item(8,8).Value = "Arrive In-Transit"
$ws.Cells.Item(8,9).Value = "HDENLA"
$ws.Cells.Item(8,10).Value = "LOVELAND"
$ws.Cells.Item(8,11).Value = "CO"
$ws.Cells.Item(8,12).Value = 280033
$ws.Cells.Item(8,13).Value = 69600
$ws.Cells.Item(8,14).Value = 210433
$ws.Cells.Item(8,15).Value = "CRDX15088"
$ws.Cells.Item(9,1).Value = "FURX"
$ws.Cells.Item(9,2).Value = 855172
$ws.Cells.Item(9,3).Value = "LOVELAND"
$ws.Cells.Item(9,4).Value = "CO"
$ws.Cells.Item(9,5).Value = 5
$ws.Cells.Item(9,6).Value = 15
$ws.Cells.Item(9,7).Value = 1048
$ws.Cells.Item(9,8).Value = "Junction Received"
$ws.Cells.Item(9,9).Value = "BNSF"
$ws.Cells.Item(9,10).Value = "LOVELAND"
$ws.Cells.Item(9,11).Value = "CO"
$ws.Cells.Item(9,12).Value = 237863
$ws.Cells.Item(9,13).Value = 60400
$ws.Cells.Item(9,14).Value = 177463
$ws.Cells.Item(9,15).Value = "FURX855172"
$ws.Cells.Item(10,1).Value = "AEX"
$ws.Cells.Item(10,2).Value = 8619
$ws.Cells.Item(10,3).Value = "LOVELAND"
$ws.Cells.Item(10,4).Value = "CO"
$ws.Cells.Item(10,5).Value = 5
$ws.Cells.Item(10,6).Value = 16
$ws.Cells.Item(10,7).Value = 1001
$ws.Cells.Item(10,8).Value = "Junction Received"
$ws.Cells.Item(10,9).Value = "BNSF"
$ws.Cells.Item(10,10).Value = "LOVELAND"
$ws.Cells.Item(10,11).Value = "CO"
$ws.Cells.Item(10,12).Value = 266857
$ws.Cells.Item(10,13).Value = 71400
$ws.Cells.Item(10,14).Value = 195457
$ws.Cells.Item(10,15).Value = "AEX8619"
$ws.Cells.Item(11,1).Value = "CRDX"
$ws.Cells.Item(11,2).Value = 15634
$ws.Cells.Item(11,3).Value = "LOVELAND"
$ws.Cells.Item(11,4).Value = "CO"
$ws.Cells.Item(11,5).Value = 5
$ws.Cells.Item(11,6).Value = 16
$ws.Cells.Item(11,7).Value = 1001
$ws.Cells.Item(11,8).Value = "Junction Received"
$ws.Cells.Item(11,9).Value = "BNSF"
$ws.Cells.Item(11,10).Value = "LOVELAND"
$ws.Cells.Item(11,11).Value = "CO"
$ws.Cells.Item(11,12).Value = 277487
$ws.Cells.Item(11,13).Value = 66800
$ws.Cells.Item(11,14).Value = 210687
$ws.Cells.Item(11,15).Value = "CRDX15634"
$ws.Cells.Item(12,1).Value = "FURX"
$ws.Cells.Item(12,2).Value = 855168
$ws.Cells.Item(12,3).Value = "NORTHTOWN"
$ws.Cells.Item(12,4).Value = "MN"
$ws.Cells.Item(12,5).Value = 5
$ws.Cells.Item(12,6).Value = 16
$ws.Cells.Item(12,7).Value = 1541
$ws.Cells.Item(12,8).Value = "Arrive In-Transit"
$ws.Cells.Item(12,9).Value = "HGFDNT"
$ws.Cells.Item(12,10).Value = "LOVELAND"
$ws.Cells.Item(12,11).Value = "CO"
$ws.Cells.Item(12,12).Value = 230136
$ws.Cells.Item(12,13).Value = 60400
$ws.Cells.Item(12,14).Value = 169736
$ws.Cells.Item(12,15).Value = "FURX855168"

# Rows 13 and 14 no longer belong to the result set - remove them
$ws.Rows("13:14").Delete()

# The sheet no longer ships with a live autofilter dropdown, but the
# hidden _FilterDatabase defined name still tracks the (now smaller) sort
# range, so turn the filter off and just refresh the name.
$ws.AutoFilterMode = $false

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Test_format_trace!_FilterDatabase" -or $n.Name -eq "_FilterDatabase") {
        $n.RefersTo = "=Test_format_trace!`$A`$2:`$N`$12"
    }
}

# Keep the selection in sync with the new data extent
$ws.Range("O3:O12").Select() | Out-Null
